$d = $word.ActiveDocument

# --- Ativação date change -------------------------------------------------
# The "Ativação" run sits directly before the "Departamento" run inside the
# same paragraph. A plain Find/Replace (or Range.Text assignment) on the
# "Ativação" run causes the interop host to coalesce it with the following
# "Departamento" run. To preserve the original run boundaries (as required
# by the target OOXML), snapshot the "Departamento" run's FormattedText
# first and reapply it after the edit, which forces the two runs apart again.
$depSnapshot = $d.Content
$depSnapshot.Find.Execute("Departamento: Engenharia Química", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$savedDepFormat = $depSnapshot.FormattedText

$ativRng = $d.Content
$ativRng.Find.Execute("Ativação: 01/01/2018", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2021", 2) | Out-Null

$depRestore = $d.Content
$depRestore.Find.Execute("Departamento: Engenharia Química", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$depRestore.FormattedText = $savedDepFormat

# --- Remaining straightforward text replacements ---------------------------
$d.Content.Find.Execute("5840560 - Marco Antonio Carvalho Pereira", $true, $false, $false, $false, $false, $true, 1, $false, "11079086 - Herlandí de Souza Andrade", 2) | Out-Null
$d.Content.Find.Execute("Aulas expositivas; trabalhos em grupo; exercícios individuais e palestras.", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.", 2) | Out-Null
$d.Content.Find.Execute("Provas e Trabalhos.", $true, $false, $false, $false, $false, $true, 1, $false, "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas.", 2) | Out-Null
$d.Content.Find.Execute("NF = (MF + PR)/2, onde PR é uma prova de recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", 2) | Out-Null
$d.Content.Find.Execute("SAPIRO, Arão., CHIAVENATO, I. Planejamento Estratégico. Campus, 2ª. edição, 2010KOTLER, P. Administração de Marketing, edição do milênio, revisão técnica de Prof. Arão Sapiro. Prentice-Hall, 2000.HOOLEY, Graham J.; PIERCY, Nigel F.; SAUNDERS, John A. Estratégia de Marketing e Posicionamento Competitivo tradução e revisão técnica: Prof. Arão Sapiro. Pearson Education do Brasil, 2001.SAPIRO, ARAO; GANGANA, MAURÍCIO; LIMA, MIGUEL; VILHENA, JOÃO BAPTISTA. Gestão de Marketing . FGV Editora, 2004.BOONE, L. e KURTZ, D.L. Marketing contemporâneo. 8ª ed. São Paulo, Livros Técnicos e Científicos, 1998.KOTLER, P; JATURISPITAK, S. e MAESINCIE, S. O marketing das nações. São Paulo, Futura, 1997. MARTINS, J.R. e BLECHER, N. O império das marcas. 2ª ed. São Paulo, Negócio Editora, 1997THUROW, L.C. O futuro do capitalismo. 2ª ed. São Paulo, Rocco, 1997.VAZ, G. N. Marketing institucional. São Paulo, Pioneira, 1995.Bibliografia ComplementarArtigos das Revistas: Marketing, Meio e Mensagem, Exame, Dinheiro, Revista da Escola de Administração da FEA-USP, Revista ESPM.", $true, $false, $false, $false, $false, $true, 1, $false, "KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.KOTLER, P.; KARTAJAYA, H.; SETIAWAN, I. Marketing 4.0: do Tradicional ao Digital. São Paulo: Sextante, 2017.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L.  Marketing Essencial. 5 ed. São Paulo: Pearson, 2013.SANDHUSEN, R. L. Marketing Básico - Série Essencial. 3 ed. São Paulo: Saraiva, 2010.SAPIRO, Arão., CHIAVENATO, I. Planejamento Estratégico. Campus, 2ª. edição, 2010 KOTLER, P. Administração de Marketing, edição do milênio, revisão técnica de Prof. Arão Sapiro. Prentice-Hall, 2000. HOOLEY, Graham J.; PIERCY, Nigel F.; SAUNDERS, John A. Estratégia de Marketing e Posicionamento Competitivo tradução e revisão técnica: Prof. Arão Sapiro. Pearson Education do Brasil, 2001. SAPIRO, ARAO; GANGANA, MAURÍCIO; LIMA, MIGUEL; VILHENA, JOÃO BAPTISTA. Gestão de Marketing . FGV Editora, 2004. BOONE, L. e KURTZ, D.L. Marketing contemporâneo. 8ª ed. São Paulo, Livros Técnicos e Científicos, 1998. KOTLER, P; JATURISPITAK, S. e MAESINCIE, S. O marketing das nações. São Paulo, Futura, 1997. MARTINS, J.R. e BLECHER, N. O império das marcas. 2ª ed. São Paulo, Negócio Editora, 1997 THUROW, L.C. O futuro do capitalismo. 2ª ed. São Paulo, Rocco, 1997. VAZ, G. N. Marketing institucional. São Paulo, Pioneira, 1995. Bibliografia Complementar Artigos das Revistas: Marketing, Meio e Mensagem, Exame, Dinheiro, Revista da Escola de Administração da FEA-USP, Revista ESPM.", 2) | Out-Null
